# Updated cryptos list on Tue May 28 19:58:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.489.72"
$ws.Range("E2").Value = "  -1.24%  "

$ws.Range("D3").Value = "'3.840.00"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'601.66"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").Value = "'170.14"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("D7").Value = "'3.840.70"
$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("E11").Value = "  +1.90%  "

$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -1.55%  "

$ws.Range("D13").Value = "'0.0000264"
$ws.Range("E13").Value = "  +4.19%  "

$ws.Range("D14").Value = "'37.19"
$ws.Range("E14").Value = "  -2.16%  "

$ws.Range("D15").Value = "'4.470.42"
$ws.Range("E15").Value = "  -1.55%  "

$ws.Range("D16").Value = "'3.836.75"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("D17").Value = "'68.411.62"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").Value = "'7.43"
$ws.Range("E19").Value = "  -1.82%  "

$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("D21").Value = "'11.15"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "'470.46"
$ws.Range("E22").Value = "  -3.56%  "

$ws.Range("D23").Value = "'0.736"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("E24").Value = "  -3.40%  "

$ws.Range("D25").Value = "'83.24"
$ws.Range("E25").Value = "  -2.30%  "

$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  -1.54%  "

$ws.Range("D27").Value = "'12.18"
$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").Value = "'2.97"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "'3.988.59"
$ws.Range("E31").Value = "  -1.18%  "

$ws.Range("D32").Value = "'7.72"
$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("D33").Value = "'31.65"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("D34").Value = "'2.32"
$ws.Range("E34").Value = "  -3.92%  "

$ws.Range("D35").Value = "'9.45"
$ws.Range("E35").Value = "  -0.32%  "

$ws.Range("D36").Value = "'3.799.47"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("E37").Value = "  -1.87%  "

$ws.Range("D38").Value = "'3.72"
$ws.Range("E38").Value = "  +14.00%  "

$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").Value = "'5.95"
$ws.Range("E41").Value = "  -2.06%  "

$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").Value = "'0.316"
$ws.Range("E43").Value = "  -3.19%  "

$ws.Range("D44").Value = "'2.00"
$ws.Range("E44").Value = "  -5.19%  "

$ws.Range("D45").Value = "'8.80"
$ws.Range("E45").Value = "  +2.02%  "

$ws.Range("B46").Value = "'Bittensor"
$ws.Range("C46").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'419.27"
$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("B47").Value = "'FLOKI"
$ws.Range("C47").Value = "'https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "'0.000294"
$ws.Range("E47").Value = "  +9.37%  "

$ws.Range("D49").Value = "'47.13"
$ws.Range("E49").Value = "  -1.90%  "

$ws.Range("D50").Value = "'26.42"
$ws.Range("E50").Value = "  +5.29%  "

$ws.Range("D51").Value = "'141.58"
$ws.Range("E51").Value = "  -1.39%  "
